# feat: add 2022-Q1 data
#
# - insert a new "2022-Q1" sheet (fund holdings for that quarter) right
#   before the "总计" (totals) sheet
# - add a new row to the top of "总计" summarizing the 2022-Q1 data
#
# Implementation notes:
#  * We duplicate the existing "总计" sheet first (so the duplicate keeps
#    identical sheetPr/pageMargins/styles), rename the ORIGINAL to
#    "2022-Q1" and reuse it for the new quarter's fund table, and rename
#    the duplicate back to "总计" for the (updated) totals table. This
#    keeps sheet order + id allocation ("2022-Q1" before "总计") correct.
#  * Numeric-looking identifiers / metrics that must stay TEXT (fund
#    code, 基金规模/股票总仓位/仓位占比/持有市值) are written with a
#    leading apostrophe (Excel's "store as text" marker) and the cell
#    Style is then reset to "Normal" so no stray number-format style is
#    left behind.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# Duplicate "总计" -> "总计 (2)", placed right after the original. This
# clone will become the refreshed "总计" sheet; the original becomes
# "2022-Q1".
$totalSheet.Copy($null, $totalSheet)
$newTotalSheet = $wb.Worksheets.Item(5)

$totalSheet.Name = "2022-Q1"
$newTotalSheet.Name = "总计"

# ---------------------------------------------------------------------
# 1) Build the "2022-Q1" fund-holdings sheet (reuses the old "总计" tab)
# ---------------------------------------------------------------------
$q1 = $totalSheet
$q1.UsedRange.Clear()

# Pull header/index formatting (style index with border + bold + center)
# from an existing fund sheet so visuals match exactly.
$styleSrc = $wb.Worksheets.Item("2021-Q4")
$styleSrc.Range("B1:H1").Copy()
$q1.Range("B1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0

$q1.Range("B2").Value = "'" + "008555"
$q1.Range("B2").Style = "Normal"

$q1.Range("C2").Value = "华商龙头优势混合"

$q1.Range("D2").Value = "'" + "1.48"
$q1.Range("D2").Style = "Normal"

$q1.Range("E2").Value = "'" + "91.13"
$q1.Range("E2").Style = "Normal"

$q1.Range("F2").Value = "'" + "5.13"
$q1.Range("F2").Style = "Normal"

$q1.Range("G2").Value = "'" + "0.0759"
$q1.Range("G2").Style = "Normal"

$q1.Range("H2").Value = 7

# ---------------------------------------------------------------------
# 2) Refresh the "总计" sheet: add 2022-Q1 as the new first data row and
#    push the previously-existing rows down by one.
# ---------------------------------------------------------------------
$tot = $newTotalSheet

$tot.Range("A2").Copy()
$tot.Range("A5").PasteSpecial(-4122)

$tot.Range("A5").Value = 3
$tot.Range("B5").Value = "2020-Q4"
$tot.Range("C5").Value = 2
$tot.Range("D5").Value = 0.01

$tot.Range("A4").Value = 2
$tot.Range("B4").Value = "2021-Q3"
$tot.Range("C4").Value = 5
$tot.Range("D4").Value = 0.57

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2021-Q4"
$tot.Range("C3").Value = 18
$tot.Range("D3").Value = 6.64

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 1
$tot.Range("D2").Value = 0.08

# Restore the originally-active sheet/tab (unrelated to this edit).
$wb.Worksheets.Item("2020-Q4").Activate()
